$wb = $excel.ActiveWorkbook

# ALC row 41 (item id 5478)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()

# ALC row 55 (item id 5517)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 420
$ws.Range("I55").Value = 312.85715
$ws.Range("J55").Value = 795
$ws.Range("K55").Value = 312.85715
$ws.Range("L55").Value = 795
$ws.Range("M55").Value = -98.85714999999999
$ws.Range("N55").Value = -1223

# ALC row 62 (item id 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2057.8333
$ws.Range("I62").Value = 1875.6
$ws.Range("J62").Value = 2969
$ws.Range("K62").Value = 1875.6
$ws.Range("L62").Value = 2969
$ws.Range("M62").Value = -1251.6
$ws.Range("N62").Value = -4217

# ALC row 65 (item id 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2057.8333
$ws.Range("I65").Value = 1875.6
$ws.Range("J65").Value = 2969
$ws.Range("K65").Value = 9378
$ws.Range("L65").Value = 14845
$ws.Range("M65").Value = -6258
$ws.Range("N65").Value = -21085

# ALC row 103 (item id 19909)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 899
$ws.Range("I103").Value = 800
$ws.Range("J103").Value = 998
$ws.Range("K103").Value = 2400
$ws.Range("L103").Value = 2994
$ws.Range("M103").Value = -1814
$ws.Range("N103").Value = -4166

# ALC row 113 (item id 27775)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2478.4167
$ws.Range("I113").Value = 1891.7142
$ws.Range("J113").Value = 3299.8
$ws.Range("K113").Value = 1891.7142
$ws.Range("L113").Value = 3299.8
$ws.Range("M113").Value = 1362.2858
$ws.Range("N113").Value = -9807.799999999999

# ALC row 118 (item id 27958)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 844.5
$ws.Range("I118").Value = 844.5
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 2533.5
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -876.5

# ALC row 132 (item id 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 681.38464
$ws.Range("I132").Value = 629.9167
$ws.Range("J132").Value = 1299
$ws.Range("K132").Value = 1889.7501
$ws.Range("L132").Value = 3897
$ws.Range("M132").Value = 640.2499
$ws.Range("N132").Value = -8957

# ALC row 138 (item id 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3920.077
$ws.Range("I138").Value = 2127.7856
$ws.Range("J138").Value = 4923.76
$ws.Range("K138").Value = 6383.3568
$ws.Range("L138").Value = 14771.28
$ws.Range("M138").Value = -1243.3568
$ws.Range("N138").Value = -25051.28

# ARM row 37 (item id 3096)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 1000
$ws.Range("I37").Value = 1000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 1000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -727

# ARM row 74 (item id 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 602.8889
$ws.Range("I74").Value = 491.29413
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 491.29413
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = 382.70587
$ws.Range("N74").Value = -4248

# ARM row 77 (item id 44000)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 602.8889
$ws.Range("I77").Value = 491.29413
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 2456.47065
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = 1911.52935
$ws.Range("N77").Value = -21236

# ARM row 102 (item id 19945)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1013
$ws.Range("I102").Value = 1013
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1013
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 609

# ARM row 132 (item id 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2050.476
$ws.Range("I132").Value = 2050.476
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6151.428
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3621.428
$ws.Range("N132").ClearContents()

# BSM row 15 (item id 1605)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 7000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 7000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 7000
$ws.Range("N15").Value = -7454

# BSM row 19 (item id 1753)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 87502.5
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 87502.5
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 87502.5
$ws.Range("N19").Value = -87848.5

# BSM row 105 (item id 19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2468.1667
$ws.Range("I105").Value = 2468.1667
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2468.1667
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -721.1667000000002

# CRP row 22 (item id 5367)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 663.36365
$ws.Range("I22").Value = 632.125
$ws.Range("J22").Value = 746.6667
$ws.Range("K22").Value = 632.125
$ws.Range("L22").Value = 746.6667
$ws.Range("M22").Value = -282.125
$ws.Range("N22").Value = -1446.6667

# CRP row 62 (item id 12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3865.8333
$ws.Range("I62").Value = 3375
$ws.Range("J62").Value = 4847.5
$ws.Range("K62").Value = 3375
$ws.Range("L62").Value = 4847.5
$ws.Range("M62").Value = -2751
$ws.Range("N62").Value = -6095.5

# CRP row 65 (item id 12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3865.8333
$ws.Range("I65").Value = 3375
$ws.Range("J65").Value = 4847.5
$ws.Range("K65").Value = 16875
$ws.Range("L65").Value = 24237.5
$ws.Range("M65").Value = -13755
$ws.Range("N65").Value = -30477.5

# CRP row 132 (item id 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1766.2858
$ws.Range("I132").Value = 1721.5
$ws.Range("J132").Value = 2348.5
$ws.Range("K132").Value = 5164.5
$ws.Range("L132").Value = 7045.5
$ws.Range("M132").Value = -2634.5
$ws.Range("N132").Value = -12105.5

# CRP row 138 (item id 42302)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 64200
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 64200
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 64200
$ws.Range("N138").Value = -74480

# CUL row 68 (item id 12895)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1030.3
$ws.Range("I68").Value = 450
$ws.Range("J68").Value = 1417.1666
$ws.Range("K68").Value = 1350
$ws.Range("L68").Value = 4251.4998
$ws.Range("M68").Value = -539
$ws.Range("N68").Value = -5873.4998

# CUL row 71 (item id 12895)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1030.3
$ws.Range("I71").Value = 450
$ws.Range("J71").Value = 1417.1666
$ws.Range("K71").Value = 4050
$ws.Range("L71").Value = 12754.4994
$ws.Range("M71").Value = 6
$ws.Range("N71").Value = -20866.4994

# CUL row 74 (item id 12859)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 20499.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 20499.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 61498.5
$ws.Range("N74").Value = -63620.5

# CUL row 77 (item id 12859)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 20499.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 20499.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 184495.5
$ws.Range("N77").Value = -195103.5

# CUL row 109 (item id 27854)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1950
$ws.Range("I109").Value = 1950
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 5850
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -4810

# CUL row 140 (item id 44097)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 6901.3335
$ws.Range("I140").Value = 1055.5
$ws.Range("J140").Value = 9824.25
$ws.Range("K140").Value = 3166.5
$ws.Range("L140").Value = 29472.75
$ws.Range("M140").Value = 2013.5
$ws.Range("N140").Value = -39832.75

# GSM row 2 (item id 5062)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.59999999999999
$ws.Range("I2").Value = 89.8
$ws.Range("J2").Value = 16.6
$ws.Range("K2").Value = 89.8
$ws.Range("L2").Value = 16.6
$ws.Range("M2").Value = 23.2
$ws.Range("N2").Value = -242.6

# GSM row 70 (item id 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5199.5
$ws.Range("I70").Value = 3700
$ws.Range("J70").Value = 6699
$ws.Range("K70").Value = 3700
$ws.Range("L70").Value = 6699
$ws.Range("M70").Value = -3430
$ws.Range("N70").Value = -7239

# GSM row 73 (item id 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5199.5
$ws.Range("I73").Value = 3700
$ws.Range("J73").Value = 6699
$ws.Range("K73").Value = 3700
$ws.Range("L73").Value = 6699
$ws.Range("M73").Value = -2764
$ws.Range("N73").Value = -8571

# GSM row 107 (item id 27802)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1629.8422
$ws.Range("I107").Value = 738.8
$ws.Range("J107").Value = 2619.889
$ws.Range("K107").Value = 738.8
$ws.Range("L107").Value = 2619.889
$ws.Range("M107").Value = 1181.2
$ws.Range("N107").Value = -6459.889

# GSM row 113 (item id 27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

# GSM row 132 (item id 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2415.5
$ws.Range("I132").Value = 2415.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7246.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4716.5

# WVR row 41 (item id 21725)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19988.25
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 19988.25
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 19988.25
$ws.Range("N41").Value = -20768.25

# WVR row 62 (item id 12589)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 23057.143
$ws.Range("I62").Value = 55002
$ws.Range("J62").Value = 10279.2
$ws.Range("K62").Value = 55002
$ws.Range("L62").Value = 10279.2
$ws.Range("M62").Value = -54378
$ws.Range("N62").Value = -11527.2

# WVR row 65 (item id 12589)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 23057.143
$ws.Range("I65").Value = 55002
$ws.Range("J65").Value = 10279.2
$ws.Range("K65").Value = 275010
$ws.Range("L65").Value = 51396
$ws.Range("M65").Value = -271890
$ws.Range("N65").Value = -57636

# WVR row 81 (item id 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2500590.2
$ws.Range("I81").Value = 786.3333
$ws.Range("J81").Value = 10000002
$ws.Range("K81").Value = 1572.6666
$ws.Range("L81").Value = 20000004
$ws.Range("M81").Value = -511.6666
$ws.Range("N81").Value = -20002126

# WVR row 84 (item id 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2500590.2
$ws.Range("I84").Value = 786.3333
$ws.Range("J84").Value = 10000002
$ws.Range("K84").Value = 7863.333000000001
$ws.Range("L84").Value = 100000020
$ws.Range("M84").Value = -2559.333000000001
$ws.Range("N84").Value = -100010628
